$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typo: "Programacion" -> "Programación" (shared string used in F8) ---
$ws.Range("F8").Value = "Programación"

# --- Update existing rows' "Ejecución" (G) and "Avance" (H) values ---
$ws.Range("G10").Value = 1.54
$ws.Range("H10").Value = 1

$ws.Range("G11").Value = 153.2
$ws.Range("H11").Value = 1

$ws.Range("G12").Value = 17
$ws.Range("H12").Value = 1

$ws.Range("G14").Value = 1515.2
$ws.Range("H14").Value = 1

# --- Prepare formatting for the 3 new rows being introduced (rows 15-19 replace old 15-16) ---
# Copy the formatting of a "numbered data" row (row 14) onto the rows that will hold
# numbered data: 15, 17, 19
$ws.Range("B14:I14").Copy()
$ws.Range("B19:I19").PasteSpecial(-4122)
$ws.Range("B17:I17").PasteSpecial(-4122)
$ws.Range("B15:I15").PasteSpecial(-4122)

# Copy the formatting of a "category" row (row 13) onto the rows that will hold
# category headers only: 16, 18
$ws.Range("B13:D13").Copy()
$ws.Range("B18:D18").PasteSpecial(-4122)
$ws.Range("B16:D16").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# --- Row 15: numbered item 5 (MR203 / Limpieza de badén) ---
$ws.Range("B15").Value = 5
$ws.Range("C15").Value = "MR203"
$ws.Range("D15").Value = "Limpieza de badén"
$ws.Range("E15").Value = "m2"
$ws.Range("F15").Value = 42.42454545454546
$ws.Range("G15").Value = 42.42
$ws.Range("H15").Value = 1
$ws.Range("I15").Value = "0+000 - 13+000"

# --- Row 16: category header (MR300 / Control de vegetación) ---
$ws.Range("B16").ClearContents()
$ws.Range("C16").Value = "MR300"
$ws.Range("D16").Value = "Control de vegetación"
$ws.Range("E16:I16").Clear()

# --- Row 17: numbered item 6 (MR301 / Roce y limpieza) ---
$ws.Range("B17").Value = 6
$ws.Range("C17").Value = "MR301"
$ws.Range("D17").Value = "Roce y limpieza"
$ws.Range("E17").Value = "m2"
$ws.Range("F17").Value = 3606.515454545454
$ws.Range("G17").Value = 3606.52
$ws.Range("H17").Value = 1
$ws.Range("I17").Value = "0+000 - 13+000"

# --- Row 18: category header (MR700 / Actividades complementarias) ---
$ws.Range("C18").Value = "MR700"
$ws.Range("D18").Value = "Actividades complementarias"

# --- Row 19: numbered item 7 (MR701 / Reparación de muros secos) ---
$ws.Range("B19").Value = 7
$ws.Range("C19").Value = "MR701"
$ws.Range("D19").Value = "Reparación de muros secos"
$ws.Range("E19").Value = "m3"
$ws.Range("F19").Value = 3.866
$ws.Range("G19").Value = 3.87
$ws.Range("H19").Value = 1
$ws.Range("I19").Value = "0+000 - 13+000"

# --- Update conditional formatting range to cover the new rows ---
$fc = $ws.Range("B7:I16").FormatConditions.Item(1)
$fc.ModifyAppliesToRange($ws.Range("B7:I19"))

# --- Update dimension reference (Excel normally manages this automatically) ---
Write-Host "Edit complete"
